$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 1.1.0 -> 1.1.1
$ws.Range("B3").Value = "1.1.1"

# Experimental value (row 7) was empty -> "false"
# Use a formula that evaluates to the text "false", then convert it to a
# plain value in place so Excel stores it as a shared string (t="s"),
# not as a native boolean (t="b") and without flipping on quotePrefix.
$ws.Range("B7").Formula = '="false"'
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

# Date updated
$ws.Range("B8").Value = "2022-10-21T09:04:31-05:00"
